# Add a new BGM resource row ("ShipMoving") for the Player & Enemy Ship
# Variety team, and leave the session on the SFX tab.

$wb  = $excel.ActiveWorkbook
$bgm = $wb.Worksheets.Item("BGM")
$sfx = $wb.Worksheets.Item("SFX")

# Insert a new row above the current row 13 (Records&Achievement /
# AchievementScreen) — this pushes the existing rows 13-15 down to 14-16
# and shifts the column-G formulas along with them.
$bgm.Rows("13:13").Insert()

# Populate the freshly inserted row 13 with the new sound-resource entry.
$bgm.Range("A13").Value = "Player&EnemyShipVariety"
$bgm.Range("B13").Value = "Ship"
$bgm.Range("C13").Value = "ShipMoving"
$bgm.Range("D13").Value = "ShipMoving.wav"
$bgm.Range("E13").Value = "S"
$bgm.Range("F13").Value = "O"

# The FileName formula for this entry hasn't been wired up yet, so it is
# left as a literal #N/A error value (matches upstream) instead of the
# usual CONCAT formula.
$bgm.Range("G13").Value = "#N/A"

# Restore the original CONCAT formula on rows 2-12 (the still-unaffected
# part of the old shared-formula block) so they stay one shared formula
# group instead of being left as per-row duplicates after the insert.
$bgm.Range("G2:G12").Formula = '=_xlfn.CONCAT("BGM_",A2,"_",B2,"_",C2,"_",D2)'

# The author ended the session on the SFX tab, not BGM.
$sfx.Activate()
